$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.164.83"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.868.47"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.60"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5061"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3915"
$ws.Range("E8").Value = "  -1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09648"
$ws.Range("E9").Value = "  -3.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").Value = "  +2.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.85"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.492"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.91"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "1.873.23"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.001"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.419"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  -1.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.92"
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06618"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.53"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.158"
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").Value = "28.228.91"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.533"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").Value = "2.089.99"
$ws.Range("E27").Value = "  +2.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.17"
$ws.Range("E28").Value = "  +1.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.77"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1057"
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.068"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.622"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.624"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.594"
$ws.Range("E35").Value = "  +4.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06732"
$ws.Range("E36").Value = "  -2.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02382"
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6344"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.973"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.175"
$ws.Range("E42").Value = "  +1.91%  "
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.55"
$ws.Range("E44").Value = "  +1.80%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6012"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.661"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.261"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.32"
$ws.Range("E48").Value = "  -0.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.992"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06830"
$ws.Range("E51").Value = "  +0.60%  "
